# Fixed typo in education categorisation for dt2
$d = $word.ActiveDocument

$d.Content.Find.Execute("7,116 (10)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "7,116 (7.9)", 2)

$d.Content.Find.Execute("22,113 (32)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "22,113 (25)", 2)

$d.Content.Find.Execute("0 (0)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "21,317 (24)", 2)

$d.Content.Find.Execute("39,697 (58)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "39,697 (44)", 2)
